$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the "Table 1." caption paragraph that used to precede the
#    table (the table now starts the document body).
# ---------------------------------------------------------------------
$cap = $d.Paragraphs(1)
$cap.Range.Delete()

# ---------------------------------------------------------------------
# Helper behaviour notes:
#  - A plain Find/Execute replace keeps a cell's run formatting intact
#    and collapses the cell text into a single run.
#  - Some cells in the target need the number text split across two
#    <w:r> runs (with identical run formatting). Toggling a character
#    property on a sub-range (set -> different value -> set back to the
#    original) forces the engine to materialise that sub-range as its
#    own run without altering the visible formatting.
# ---------------------------------------------------------------------

function Replace-Text($old, $new) {
    $r = $d.Content
    $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

function Split-RunAt($text, $offset) {
    # Re-find the (already replaced) text and force a run split after
    # $offset characters by toggling the font size of that leading
    # sub-range back and forth.
    $r = $d.Content
    $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $orig = $r.Font.Size
    $splitPoint = $r.Start + $offset
    $sub = $d.Range($r.Start, $splitPoint)
    $sub.Font.Size = $orig + 1
    $sub.Font.Size = $orig
}

# ---------------------------------------------------------------------
# 2. Numeric corrections inside the table.
# ---------------------------------------------------------------------

# 165.9 -> 166   (split into "16" + "6")
# Do this before the 165.6 -> 166 replacement below so "166" is still
# unique in the document when we search for it to apply the split.
Replace-Text "165.9" "166"
Split-RunAt "166" 2

# 165.6 -> 166   (single run)
Replace-Text "165.6" "166"

# 171.8 -> 172   (split into "17" + "2")
Replace-Text "171.8" "172"
Split-RunAt "172" 2

# 178.1 -> 178   (single run)
Replace-Text "178.1" "178"

# 186.4 -> 186   (single run)
Replace-Text "186.4" "186"

# minus sign below is a literal U+2212 (MINUS SIGN), not a hyphen

# -13.0 -> -13   (single run)
Replace-Text "−13.0" "−13"

# -13.3 -> -13   (single run)
Replace-Text "−13.3" "−13"

# -12.4 -> -12.5 (split into "-12." + "5")
Replace-Text "−12.4" "−12.5"
Split-RunAt "−12.5" 4

# -12.1 -> -12   (single run)
Replace-Text "−12.1" "−12"
